# Weekly fruit/vegetable data refresh: a new week's worth of records is
# inserted for this Cilantro subset (Vega Central Mapocho de Santiago),
# pushing the existing rows down by two rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right before row 926; this shifts the existing
# rows 926:1001 down to 928:1003 and extends the sheet dimension
# accordingly (matches the A1:R1001 -> A1:R1003 dimension change).
$ws.Range("A926:A927").EntireRow.Insert()

# --- New row 926 ---
$ws.Cells.Item(926, 1).Value = 9
$ws.Cells.Item(926, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(926, 3).Value = "Metropolitana"
$ws.Cells.Item(926, 4).Value = 45106
$ws.Cells.Item(926, 5).Value = 13
$ws.Cells.Item(926, 6).Value = 100112040
$ws.Cells.Item(926, 7).Value = "Cilantro"
$ws.Cells.Item(926, 8).Value = "Sin especificar"
$ws.Cells.Item(926, 9).Value = "Primera"
$ws.Cells.Item(926, 10).Value = 52
$ws.Cells.Item(926, 11).Value = 8000
$ws.Cells.Item(926, 12).Value = 8000
$ws.Cells.Item(926, 13).Value = 8000
$ws.Cells.Item(926, 14).Value = "`$/caja 36 atados"
$ws.Cells.Item(926, 15).Value = "Región Metropolitana"
$ws.Cells.Item(926, 16).Value = 222
$ws.Cells.Item(926, 17).Value = 36
$ws.Cells.Item(926, 18).Value = "Hortaliza"

# --- New row 927 ---
$ws.Cells.Item(927, 1).Value = 9
$ws.Cells.Item(927, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(927, 3).Value = "Metropolitana"
$ws.Cells.Item(927, 4).Value = 45106
$ws.Cells.Item(927, 5).Value = 13
$ws.Cells.Item(927, 6).Value = 100112040
$ws.Cells.Item(927, 7).Value = "Cilantro"
$ws.Cells.Item(927, 8).Value = "Sin especificar"
$ws.Cells.Item(927, 9).Value = "Primera"
$ws.Cells.Item(927, 10).Value = 160
$ws.Cells.Item(927, 11).Value = 16000
$ws.Cells.Item(927, 12).Value = 18000
$ws.Cells.Item(927, 13).Value = 17000
$ws.Cells.Item(927, 14).Value = "`$/docena de atados"
$ws.Cells.Item(927, 15).Value = "Región Metropolitana"
$ws.Cells.Item(927, 16).Value = 5667
$ws.Cells.Item(927, 17).Value = 3
$ws.Cells.Item(927, 18).Value = "Hortaliza"

# Make sure the new date cells carry the same date style as the rest of
# column D (style index 2 in the original workbook / "yyyy-mm-dd" style).
$ws.Range("D926:D927").NumberFormat = $ws.Range("D928").NumberFormat
